$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 74 (Serie 01-01-2021) values for columns B, C, D, F
$ws.Range("B74").Value = -919
$ws.Range("C74").Value = -507
$ws.Range("D74").Value = 1445
$ws.Range("F74").Value = -2639

# Add new row 75 (Serie 01-04-2021)
# Enter the date-looking label via a text formula so Excel does not
# auto-convert it to a date serial, then flatten it to a plain value
# with Paste Special so no new number-format style gets introduced.
$ws.Range("A75").Formula = "=""01-04-2021"""
$ws.Range("A75").Copy()
$ws.Range("A75").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("B75").Value = -4795
$ws.Range("C75").Value = -3436
$ws.Range("D75").Value = -11714
$ws.Range("E75").Value = 74
$ws.Range("F75").Value = 5835
$ws.Range("G75").Value = 4446
